$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Comment.Text('Op te leveren producten (M01)

Projecten zijn gesplitst in een voorbereidingsfase en een realisatiefase (zie [M14: Projecten splitsen in een voorbereidingsfase en een realisatiefase](#projecten-splitsen-in-een-voorbereidingsfase-en-een-realisatiefase-m14-)); alle fasen kennen specifieke producten. De onderstaande tabel benoemt de producten die een project in die fasen moet realiseren of waarover het project moet kunnen beschikken, indien er andere auteurs zijn (bijvoorbeeld de opdrachtgever). Het project zorgt voor een volledige oplevering van alle genoemde producten, waaronder de producten die niet door het project zelf zijn gemaakt.

Als tijdens een project bestaande software dient te worden afgebouwd, onderhouden en/of herbouwd, vindt een onderzoek plaats naar de compleetheid en consistentie van de bestaande softwareproducten aan de hand van de onderstaande tabel (inclusief de deliverables in de kolom ''Realisatiefase'') en wordt de kwaliteit van de bestaande software-producten getoetst ([M02: Continu voldoen aan kwaliteitsnormen](#continu-voldoen-aan-kwaliteitsnormen-m02-)). Dit onderzoek is onderdeel van de voorbereidingsfase en wordt uitgevoerd door vertegenwoordigers van de projectenorganisatie en medewerkers van het desbetreffende project.

| Product                | Voorbereidings᠆fase  | Voorbereidings᠆fase met onderzoek  | Realisatie᠆fase |
|----------------------------------------------------------------------------------------|---|---|---|
| BIA (Business impact analysis)*                                                        | ✔ | ✔ | ✔ |
| PIA (Privacy impact analysis)**                                                        | ✔ | ✔ | ✔ |
| Beschrijving van functionele eisen                                                     | ✔ | ✔ | ✔ |
| Beschrijving van niet-functionele eisen                                                | ✔ | ✔ | ✔ |
| Ontwerp- en architectuurdocumentatie (software, interactie, infrastructuur)            | ✔ | ✔ | ✔ |
| Testdocumentatie: testplannen                                                          | ✔ | ✔ | ✔ |
| Testdocumentatie: testgevallen, rapportages                                            |   |   | ✔ |
| Informatiebeveiligingsplan                                                             | ✔ | ✔ | ✔ |
| Projectplan                                                                            | ✔ | ✔ |   |
| Kwaliteitsplan                                                                         | ✔ | ✔ |   |
| Deploybare versie van de software                                                      |   |   | ✔ |
| Broncode, inclusief de benodigdheden voor het bouwen van de software                   |   |   | ✔ |
| Regressietests, inclusief de benodigdheden voor het uitvoeren van de regressietesten   |   |   | ✔ |
| Vrijgaveadvies                                                                         |   |   | ✔ |
| Release notes                                                                          |   |   | ✔ |
| Deploymentdocumentatie                                                                 |   |   | ✔ |
| Uitkomsten onderzoek (bevindingen, risico''s, mitigerende maatregelen)                  |   | ✔ |   |
| Transitieplan voor af te bouwen, te onderhouden en/of te herbouwen softwareproducten   |   | ✔ |   |
| Plan voor aflossen technische schuld, indien van toepassing                            |   | ✔ |   ||

*) De opdrachtgever stelt bij voorkeur een BIA en PIA voorafgaand aan een voorbereidende fase beschikbaar. Indien dat niet mogelijk is worden deze analyses gedurende de voorbereidende fase onder verantwoordelijkheid van de opdrachtgever uitgevoerd.
**) Indien een PIA niet nodig is, is een verklaring daaromtrent vereist.

Rationale

Het uniformeren van op te leveren producten biedt voordelen voor planning (het is bekend welke producten gemaakt moeten worden), voor bemensing (het is bekend welke expertise nodig is) en voor het uitwisselen van medewerkers.

De voorgeschreven producten stellen de ontvanger in staat om de opgeleverde software uit te voeren, te beheren en te onderhouden. Daarnaast is duidelijk welke eventueel openstaande punten er nog zijn. De voorgeschreven producten bieden voldoende verantwoording richting de ontvanger voor uitgevoerde werkzaamheden.

De genoemde producten uit de voorbereidingsfase hebben tot doel om enerzijds de omvang, kosten en doorlooptijd van de realisatiefase te kunnen schatten en anderzijds om de kaders voor de realisatiefase te bepalen, zodat de scope, aanpak en oplossingsrichting in grote lijnen bekend zijn.

Een BIA en eventuele PIA zijn richtinggevend voor de in de voorbereidingsfase te selecteren beveiligingsmaatregelen en zijn daarom, bij voorkeur, voorafgaand aan het project al beschikbaar.

In een BIA legt de vragende organisatie vast hoe belangrijk informatiebeveiliging is voor de eigen bedrijfsvoering/processen. Naast de gevoeligheid voor incidenten komt hierin ook de ''risk appetite'' van de organisatie tot uiting. Alleen de organisatie zelf kan hierover een uitspraak doen.

In een PIA legt de vragende organisatie vast wat de privacy-gevoeligheid is van de gegevens die in een proces of systeem worden verzameld en verwerkt. Zicht op privacygevoelige gegevens en het (laten) treffen van adequate en afdoende beschermingsmaatregelen is een wettelijke plicht die een organisatie niet aan een andere partij kan verdragen.

ICTU

ICTU hanteert de volgende documenten, templates en documentstandaarden voor softwarerealisatieprojecten:

- De beschrijving van niet-functionele eisen is gebaseerd op ISO (International Organization for Standardization)-25010, de BIR (Baseline Informatiebeveiliging Rijksdienst) en de methode Grip op SSD (Secure software development) van het CIP (Centrum Informatiebeveiliging en Privacybescherming), en bevat een prioritering van de niet-functionele eisen. De beschrijving van niet-functionele eisen is gebaseerd op het ICTU NFE (Niet-functionele eisen)-template. De beschrijving bevat in ieder geval eisen aan toegangsbeveiliging, aan beheerfuncties, aan logging en aan het gewenste gedrag van de software bij uitval van infrastructurele diensten, zoals een log-server;

- De beschrijving van functionele eisen bestaat uit een geprioriteerde backlog met epics en/of user stories. De beschrijving bevat in ieder geval eisen voor (ondersteuning van) beheerfuncties die door de beoogd beheerder gesteld worden en voor logging, inclusief de (globale) inhoud van te loggen business events (gebeurtenissen op procesniveau) en de daarvoor geldende bewaartermijnen;

- De ontwerp- en architectuurdocumentatie bestaat uit een PSA (Projectstartarchitectuur), een SAD (Softwarearchitectuurdocument), een IA (Infrastructuurarchitectuur), een GFO (Globaal functioneel ontwerp) bijvoorbeeld in de vorm van use cases, en een prototype en/of interactieontwerp. De architectuurdocumenten moeten expliciet inzichtelijk maken hoe aan de niet-functionele eisen wordt voldaan door uit te werken welke oplossingen en mechanieken gekozen zijn, bijvoorbeeld voor identificatie, authenticatie, autorisatie, concurrency, transactionele verwerking of logging;

- De testdocumentatie bestaat uit een mastertestplan, gemaakt op basis van een PRA (Productrisicoanalyse). Beveiligingstesten zijn een integraal onderdeel van het mastertestplan en worden als zodanig afgestemd met de opdrachtgever;

- Het informatiebeveiligingsplan is gebaseerd op een dreigingen- en kwetsbaarhedenanalyse (TVA (Threat and vulnerability assessment)) en bevat een maatregelenselectie informatiebeveiliging. De TVA wordt tijdens de voorfase opgesteld op basis van de resultaten van de BIA, de eventuele PIA en inhoud van de ontwerp- en architectuurdocumentatie. Een TVA levert een deel van een traceerbare onderbouwing voor de te treffen beveiligingsmaatregelen;

- Het vrijgaveadvies bevat ten minste alle nog openstaande testbevindingen en geconstateerde beveiligingsbevindingen. Zie ook [M26: Periodieke beoordeling informatiebeveiliging](#periodieke-beoordeling-informatiebeveiliging-m26-) en [M16: Verplichte tools](#verplichte-tools-m16-). Indien er beveiligingsissues zijn, zijn deze voorzien van een beschreven voorziene impact.

- De deploymentdocumentatie bevat informatie over de eisen die een applicatie stelt aan een omgeving en de stappen die nodig zijn om de applicatie in die omgeving veilig te installeren en configureren. De documentatie bevat daartoe onder meer aanwijzingen voor de HTTP-header en -request-configuratie van de webserver en voor het verwijderen van overbodige header-informatie zoals de ''Server''-header. Ook zijn er aanwijzingen voor veilige configuratie(s) van (externe) toegang tot de beheerinterface. De documentatie bevat daarnaast in ieder geval een beschrijving van de protocollen en services die de applicatie aanbiedt, de protocollen, services en accounts die het product gebruikt en de protocollen, services en accounts die de applicatie gebruikt voor beheer;

Zie de bijlage [Documenten voor M01: Op te leveren producten](#documenten-voor-m01-op-te-leveren-producten) voor een uitgebreider overzicht van de documenten en documentstandaarden die ICTU hanteert voor softwarerealisatieprojecten.

Het genoemde onderzoek voert ICTU uit als onderdeel van een "due diligence". Een due diligence wordt uitgevoerd in samenwerking met een potentiële opdrachtgever en biedt, naast het genoemde onderzoek, ook de opdrachtgever de kans zich een oordeel te vormen over de werkwijze van ICTU en de verwachte samenwerking.
')
$ws.Range("B24").Comment.Text('Gebruik van ISO-25010 (M13)

De standaard ISO/IEC (International Electrotechnical Commission) 25010:2011, kortweg "ISO-25010", biedt een model voor het beschrijven van productkwaliteit. Kwaliteitseigenschappen zijn voorzien van een naam, definitie en classificatie. ISO-25010 dekt een breed spectrum van kwaliteitseigenschappen af.

Voor specificatie en documentatie van vereiste en gewenste kwaliteitseigenschappen, de niet-functionele eisen, maken projecten gebruik van de terminologie uit ISO-25010. Projecten gebruiken ISO-25010 om te controleren of alle relevante kwaliteitseigenschappen van het op te leveren eindproduct worden meegenomen in de ontwikkeling en/of onderhoud van het product.

Rationale

De standaard ISO-25010 biedt een model voor productkwaliteit. De standaard biedt geen concrete maatregelen, maar biedt wel een begrippenkader en dekt het volledige spectrum van mogelijk relevante kwaliteitseigenschappen af. Het gebruiken van een standaard voor specificatie van kwaliteit voorkomt miscommunicatie over kwaliteitseigenschappen en de breedte van de standaard zorgt ervoor dat alle relevante aspecten aan bod komen.

ICTU

ICTU gebruikt ISO-25010 voor documentatie en specificatie van productkwaliteit.
')
$ws.Range("B25").Comment.Text('Continu voldoen aan kwaliteitsnormen (M02)

Producten voldoen zo snel mogelijk vanaf de start van een project aan de door het project en projectenorganisatie vastgestelde kwaliteitsnormen. De kwaliteit van producten die nog niet zijn afgerond of nog niet aan de normen voldoen, wordt door het project bewaakt. Herstel van de kwaliteit wordt planmatig opgepakt (zie ook [M08: Technische schuld](#technische-schuld-m08-)).

De kwaliteitsnormen zijn in deze versie van de kwaliteitsaanpak nog niet opgenomen.

Rationale

Het zo snel mogelijk en continu voldoen aan de kwaliteitsnormen beperkt toekomstige hersteltijd. Het dwingt tevens een structurele kwaliteitscontrole af.

ICTU

Bij ICTU wordt tijdens de voorfase van softwarerealisatieprojecten (zie ook [M14: Projecten splitsen in een voorbereidingsfase en een realisatiefase](#projecten-splitsen-in-een-voorbereidingsfase-en-een-realisatiefase-m14-)) het voldoen aan de kwaliteitsnormen met behulp van reviews gecontroleerd. Tijdens de realisatiefase van softwarerealisatieprojecten wordt het voldoen aan de kwaliteitsnormen diverse malen per uur gemeten door het ''Kwaliteitssysteem'', genaamd HQ (Holistic Software Quality Reporting). Het project kijkt dagelijks of er afwijkingen van de normen zijn en onderneemt actie indien nodig. Ook de kwaliteitsmanager signaleert afwijkingen en meldt deze bij het project. De ICTU-specifieke invulling van de kwaliteitsnormen is te vinden in het helpmenu van de geautomatiseerde kwaliteitsrapportages van ICTU.
')
$ws.Range("B28").Comment.Text('Periodieke beoordeling informatiebeveiliging (M26)

Projecten laten periodiek de beveiliging van de ontwikkelde software beoordelen. De code wordt zowel geautomatiseerd als handmatig door een beveiligingsexpert onderzocht op veelvoorkomende kwetsbaarheden en op het voldoen aan voorgeschreven normen. Overheidsspecifieke beveiligingsnormen of -raamwerken, zoals de BIR (Baseline Informatiebeveiliging Rijksdienst), bieden een basis voor de beoordeling.
De projectorganisatie zorgt ervoor dat deze expertise op afroep beschikbaar gesteld kan worden aan projecten. Bevindingen uit de beveiligingstest worden vastgelegd als onderdeel van de werkvoorraad voor het ontwikkelproces (zie [M05: Iteratief en incrementeel ontwikkelproces](#iteratief-en-incrementeel-ontwikkelproces-m05-)).

Rationale

Door het inschakelen van actuele, specifieke expertise wordt de kans vergroot dat eventuele kwetsbaarheden in de gerealiseerde software tijdig herkend worden.

ICTU

Software wordt minimaal bij iedere grote release of ten minste twee keer per jaar onderworpen aan een beveiligingstest door beveiligingsexperts die ICTU daarvoor inhuurt. Op basis van documentatie en architectuurstudie, crystalbox security audits (broncodescan) en penetratieaudits beoordelen deze experts of de software voldoet aan de projectspecifieke niet-functionele eisen met betrekking tot beveiliging, of bekende kwetsbaarheden (zoals bijvoorbeeld in de OWASP Top 10 genoemd) vermeden zijn en in hoeverre voldoende invulling gegeven is aan de normen vanuit die vanuit BIR en SSD gelden.

Indien door de opdrachtgever gewenst kunnen securitytesten door een onafhankelijke derde partij worden uitgevoerd in een daarvoor door de opdrachtgever beschikbaar gestelde omgeving. Dit kan zowel incidenteel als structureel worden ingericht. Afspraken hierover worden bij voorkeur al in de voorbereidingsfase gemaakt.

De beveiligingstesten vinden plaats in aanvulling op de door tools uitgevoerde continue beveiligingsanalyse van de gerealiseerde software, zie [M16: Verplichte tools](#verplichte-tools-m16-). Bevindingen uit zowel een beveiligingstest als de continue analyse worden in Jira als issue - gemarkeerd als beveiligingsbugreport - vastgelegd op de backlog van het project.
')
$ws.Range("B36").Comment.Text('Continuous delivery pipeline (M07)

Er is een geautomatiseerde continuous delivery pipeline die aantoonbaar correct werkt en ten minste de volgende activiteiten uitvoert:

- bouw van de software,
- unit tests,
- regressietests,
- kwaliteitscontroles,
- performancetests (*),
- beveiligingstests (*),
- installatie van de software,
- oplevering van het totale product, dus inclusief alle deliverables, in de vorm zoals bruikbaar voor en afgesproken met de opdrachtgever.

(*) Idealerwijs voert de geautomatiseerde continuous delivery pipeline ook performance tests en beveiligingstests uit. Vanwege de doorlooptijden van tests (met name van duurtesten) en licenties van testtools is dat niet altijd haalbaar. In dat geval vinden de performance tests en beveiligingstests periodiek en zo vaak mogelijk plaats, bij voorkeur dagelijks.

De projectenorganisatie voorziet in mensen en hulpmiddelen, zodat projecten deze pipeline kunnen toepassen. Projecten zijn verantwoordelijk voor de correcte werking van de pipelin.

Rationale

Software incrementeel opleveren (zie [M05: Iteratief en incrementeel ontwikkelproces](#iteratief-en-incrementeel-ontwikkelproces-m05-)) vereist dat de software frequent gebouwd, getest en opgeleverd kan worden. Om dit efficiënt en foutvrij te doen, dient het proces van bouwen, testen en opleveren geautomatiseerd te zijn; een continuous delivery pipeline faciliteert dit.

ICTU

ICTU gebruikt Jenkins of TFS (Team Foundation Server) als tool voor de implementatie van de continuous delivery pipeline. De ICTU Release Manager ondersteunt de laatste stap (oplevering van het totale product).
')
$ws.Range("B45").Comment.Text('Technische schuld (M08)

Technische schuld is inzichtelijk en wordt planmatig aangepakt. De kwaliteitsmanager is verantwoordelijk voor het inzichtelijk maken van de technische schuld. De projectverantwoordelijke is verantwoordelijk voor het planmatig aanpakken van de technische schuld.

Rationale

De aanwezigheid van technische schuld heeft nadelige invloed op de kwaliteit van de eindproducten. Anderzijds is het ontstaan van technische schuld gedurende een project vaak onvermijdelijk. Het is daarnaast ook mogelijk dat een deel van de technische schuld bij aanvang van het project al bestond en mogelijk niet wordt opgelost. In alle gevallen is het verstandig om te weten welke technische schuld bestaat. Om te voorkomen dat technische schuld niet wordt opgelost en uitsluitend toeneemt, is het zaak om het verminderen van technische schuld planmatig aan te pakken.

ICTU

ICTU gebruikt HQ (Holistic Software Quality Reporting), een door ICTU ontwikkeld, open source, geautomatiseerd kwaliteitssysteem, om bestaande technische schuld inzichtelijk te maken en de planning van het aflossen van de schuld vast te leggen, voor zover het technische schuld betreft van kwaliteitseigenschappen die HQ kan meten.
')
$ws.Range("B55").Comment.Text('Open source tools (M15)

Bij de selectie van tools ter ondersteuning van de projectuitvoering  geeft de projectenorganisatie voorkeur aan open source tools.

Rationale

Conform de rationale uit NORA (Nederlandse Overheid Referentiearchitectuur) voor het gebruik van open source tools, zoals beschreven in NORA v3.0 drijfveer “beleid open standaarden” ([http://www.noraonline.nl/wiki/Beleid_open_standaarden](http://www.noraonline.nl/wiki/Beleid_open_standaarden)).

ICTU

Tools die ICTU ontwikkelt ter ondersteuning van softwarerealisatieprojecten, worden bij voorkeur als open source beschikbaar gesteld.
')
$ws.Range("B56").Comment.Text('Verplichte tools (M16)

De projectenorganisatie stelt het gebruik van de volgende tools voor alle projecten verplicht:

1. Een tool dat agile werken ondersteunt. Een dergelijk tool voorziet in het opvoeren van eisen, het opvoeren van logische testgevallen en het koppelen van logische testgevallen aan eisen, het bijhouden van een werkvoorraad, het plannen van iteraties en het toewijzen van eisen aan iteraties.
2. Een tool dat het inrichten en uitvoeren van een continuous delivery pipeline ondersteunt.
3. Een tool dat het monitoren van de kwaliteit van broncode ondersteunt.
4. Een tool dat het releasen van software ondersteunt.
5. Een tool dat het maken van testrapportages ondersteunt.
6. Een tool dat het maken van kwaliteitsrapportages ondersteunt.
7. Een tool dat de configuratie van de applicatie en de omgeving, waarbinnen die applicatie draait, controleert op bekende en veelvoorkomende kwetsbaarheden.
8. Een tool dat de door de applicatie gebruikte versies van externe bibliotheken, raamwerken of andersoortige bouwblokken scant op bekende kwetsbaarheden.
9. Een tool dat de broncode geautomatiseerd beoordeelt op het voorkomen van bekende kwetsbare constructies.

Rationale

Projecten hebben een redelijke vrijheid bij het kiezen van tools, maar het gebruik van een aantal is verplicht gesteld. Deze tools zijn nodig voor een efficiënte uitvoering van deze kwaliteitsaanpak. Uniform gebruik van deze tools maakt het mogelijk koppeling tussen die tools voor alle projecten te standaardiseren. Daarnaast bevordert het de uitwisselbaarheid van medewerkers en neemt het risico op het gebruik van onvolwassen tools af.

ICTU

ICTU gebruikt hiervoor de volgende tools:

1. Jira - De ''eisen'' worden, conform Scrumterminologie, geregistreerd als epics en/of user stories, de werkvoorraad als backlog, de iteraties als sprints,
2. Jenkins voor Javaprojecten en TFS (Team Foundation Server) voor DotNet-projecten,
3. SonarQube,
4. Release Manager,
5. Reporting (ontwikkeld met behulp van BIRT (Business Intelligence Reporting Tool)),
6. Kwaliteitssysteem (HQ),
7. OpenVAS (Vulnerability Assessment System) en OWASP (Open Web Application Security Project) ZAP (Zed Attack Proxy),
8. OWASP Dependency Checker,
9. Checkmarx.
')
$ws.Range("B71").Comment.Text('Afgeschermde digitale omgeving (M19)

De projectenorganisatie geeft de projecten de beschikking over eigen, afgeschermde digitale omgevingen, waarbinnen ze de door het project ontwikkelde software en tools kunnen installeren.

Rationale

Door het bieden van een afgeschermde digitale omgeving zijn de afhankelijkheden en invloeden tussen projecten minimaal en worden beveiligingsrisico''s verkleind.

ICTU

ICTU ondersteunt dit met Docker en/of virtuele machines en een VLAN (Virtual local area network) per project. Een nieuwe afgeschermde digitale omgeving is binnen een werkweek na aanvraag beschikbaar.
')

$ws.Range("A1").Value = 'Onderstaande checklist kan gebruikt worden voor het uitvoeren van een assessment tegen de Kwaliteitsaanpak ICTU Software Realisatie versie 1.2.1-build.38, 31-08-2018.'
$ws.Range("B5").Value = '1. BIA (Business impact analysis)*'
$ws.Range("B6").Value = '2. PIA (Privacy impact analysis)**'
$ws.Range("B58").Value = '2. Een tool dat het inrichten en uitvoeren van een continuous delivery pipeline ondersteunt. Bij ICTU: Jenkins voor Javaprojecten en TFS (Team Foundation Server) voor DotNet-projecten,'
$ws.Range("B61").Value = '5. Een tool dat het maken van testrapportages ondersteunt. Bij ICTU: Reporting (ontwikkeld met behulp van BIRT (Business Intelligence Reporting Tool)),'
$ws.Range("B63").Value = '7. Een tool dat de configuratie van de applicatie en de omgeving, waarbinnen die applicatie draait, controleert op bekende en veelvoorkomende kwetsbaarheden. Bij ICTU: OpenVAS (Vulnerability Assessment System) en OWASP (Open Web Application Security Project) ZAP (Zed Attack Proxy),'
